# Auto-generated edit script: update cryptos list (Sun Dec 17 02:59:36 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.993.48"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.218.08"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.29"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.52%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.69"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.83%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.613"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.33"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0959"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.05%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.550.50"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.26"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.843"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.232.24"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.901.62"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000110"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +12.16%  "
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.19"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.24"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.66%  "
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.18"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +32.14%  "
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "229.03"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.12"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -8.16%  "
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.58"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.20%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "WEMIXToken"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.59"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.19%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.28"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.84%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.96"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.53%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.61"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.60"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +14.23%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0798"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.22%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.125"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "29.47"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.14%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.114"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.99%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.31"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.30%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0301"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.03"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.75%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.14"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.75%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "64.60"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.18%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.65"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.08%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.200"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.70"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.43%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.59"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.61%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.101"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.84%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.40"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +6.20%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.12"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.16"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.71"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.75%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.427.14"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.96%  "
